# Adds the portfolio-page link to the contact-info block of the resume.
#
# Before:
#   ... Email: ... paragraph (ContactInfoEmphasis, pPr already has rPr sz/szCs 20)
#   " LinkedIn: https://www.linkedin.com/in/cathy-colette-tanya-32b10b217/  " paragraph
#       (ContactInfoEmphasis, pPr has NO rPr yet) + trailing divider-dot sdt
#   </w:tc>  (end of the contact-info table cell)
#
# After:
#   ... Email: ... paragraph (unchanged)
#   " LinkedIn: ... " paragraph, now with <w:rPr><w:sz 20/><w:szCs 20/></w:rPr> in its pPr
#   new paragraph: "https://cathytanya.github.io/cathy-colette-tanya-portfolio/"
#       (ContactInfoEmphasis, pPr has rPr sz/szCs 20, run has matching rPr)
#   new empty paragraph (ContactInfoEmphasis, pPr with contextualSpacing only)
#   </w:tc>

$d = $word.ActiveDocument

# --- locate the LinkedIn paragraph -----------------------------------------
# It is the paragraph whose text contains "LinkedIn:" inside the contact-info
# table at the top of the document.
$linkedInIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*LinkedIn:*") {
        $linkedInIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($linkedInIndex)
$pStart = $p.Range.Start
$pEnd = $p.Range.End

# --- Step 1: rewrite the LinkedIn paragraph, adding the paragraph-mark ------
# run properties (w:pPr/w:rPr sz=20/szCs=20) while leaving every run exactly
# as it was (same text, same rsid attributes) - only the pPr gains the rPr.
$rng = $d.Range($pStart, $pEnd)

$linkedInXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml"><w:body><w:p w14:paraId="7569A12C" w14:textId="39943CAA" w:rsidR="00692703" w:rsidRPr="0077494A" w:rsidRDefault="00692703" w:rsidP="00913946"><w:pPr><w:pStyle w:val="ContactInfoEmphasis"/><w:contextualSpacing w:val="0"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="0077494A"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0077494A" w:rsidRPr="0077494A"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">LinkedIn: https://www.linkedin.com/in/cathy-colette-tanya-32b10b217/ </w:t></w:r><w:r w:rsidRPr="0077494A"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:sdt><w:sdtPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:alias w:val="Divider dot:"/><w:tag w:val="Divider dot:"/><w:id w:val="759871761"/><w:placeholder><w:docPart w:val="5A939148F715498DA268CF33180BE177"/></w:placeholder><w:temporary/><w:showingPlcHdr/><w15:appearance w15:val="hidden"/></w:sdtPr><w:sdtContent><w:r w:rsidRPr="0077494A"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>&#183;</w:t></w:r></w:sdtContent></w:sdt></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($linkedInXml)

# --- Step 2: insert the new portfolio-link paragraph and the trailing blank
# paragraph right after the (just rewritten) LinkedIn paragraph. Re-fetch the
# paragraph fresh from the collection since the previous InsertXML was a
# structural edit and can invalidate old handles.
$p2 = $d.Paragraphs.Item($linkedInIndex)
$afterRng = $d.Range($p2.Range.End, $p2.Range.End)

$newParasXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ContactInfoEmphasis"/><w:contextualSpacing w:val="0"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>https://cathytanya.github.io/cathy-colette-tanya-portfolio/</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ContactInfoEmphasis"/><w:contextualSpacing w:val="0"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$afterRng.InsertXML($newParasXml)
